$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the obsolete "(census results)" caption in A2 and the empty
# formatted cell next to it, as well as the stray formatted cell in B1.
$ws.Range("A2:B2").Clear()
$ws.Range("B1").Clear()

# Remove the now-empty spacer row (old row 3) that separated the title
# block from the data table.
$ws.Rows(3).Delete()

# Drop the 1989 and 2002 columns, keeping only the 2014 figures.
$ws.Columns("B:C").Delete()

# Rename the sheet to the municipality name.
$ws.Name = "წალკა"

$ws.Range("A2").Select() | Out-Null
